# Apply updated cryptocurrency price/volume data to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference style (default/normal, no borders or bold) taken from an
# untouched data cell, used to keep text-forced numeric-looking cells
# from picking up a "quote prefix" style when we set them as text.
$normalStyle = $ws.Range("B2").Style

$ws.Range("D2").Value = "'65.514.51"
$ws.Range("D2").Style = $normalStyle
$ws.Range("E2").Value = "  -1.98%  "

$ws.Range("D3").Value = "'3.384.69"
$ws.Range("D3").Style = $normalStyle
$ws.Range("E3").Value = "  -2.65%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").Value = "'595.68"
$ws.Range("D5").Style = $normalStyle
$ws.Range("E5").Value = "  -1.48%  "

$ws.Range("D6").Value = "'141.31"
$ws.Range("D6").Style = $normalStyle
$ws.Range("E6").Value = "  -4.92%  "

$ws.Range("D7").Value = "'1.00"
$ws.Range("D7").Style = $normalStyle
$ws.Range("E7").Value = "  -0.07%  "

$ws.Range("D8").Value = "'3.382.06"
$ws.Range("D8").Style = $normalStyle
$ws.Range("E8").Value = "  -2.70%  "

$ws.Range("E9").Value = "  -3.19%  "

$ws.Range("D10").Value = "'7.91"
$ws.Range("D10").Style = $normalStyle
$ws.Range("E10").Value = "  +4.94%  "

$ws.Range("E11").Value = "  -6.94%  "

$ws.Range("D12").Value = "'0.404"
$ws.Range("D12").Style = $normalStyle
$ws.Range("E12").Value = "  -5.05%  "

$ws.Range("D13").Value = "'3.961.39"
$ws.Range("D13").Style = $normalStyle
$ws.Range("E13").Value = "  -2.53%  "

$ws.Range("D14").Value = "'0.0000198"
$ws.Range("D14").Style = $normalStyle
$ws.Range("E14").Value = "  -7.68%  "

$ws.Range("D15").Value = "'29.53"
$ws.Range("D15").Style = $normalStyle
$ws.Range("E15").Value = "  -7.17%  "

$ws.Range("E16").Value = "  -0.61%  "

$ws.Range("D17").Value = "'65.475.42"
$ws.Range("D17").Style = $normalStyle
$ws.Range("E17").Value = "  -2.13%  "

$ws.Range("D18").Value = "'3.384.05"
$ws.Range("D18").Style = $normalStyle
$ws.Range("E18").Value = "  -2.65%  "

$ws.Range("D19").Value = "'10.31"
$ws.Range("D19").Style = $normalStyle
$ws.Range("E19").Value = "  +1.82%  "

$ws.Range("D20").Value = "'6.08"
$ws.Range("D20").Style = $normalStyle
$ws.Range("E20").Value = "  -6.09%  "

$ws.Range("D21").Value = "'14.54"
$ws.Range("D21").Style = $normalStyle
$ws.Range("E21").Value = "  -5.91%  "

$ws.Range("D22").Value = "'412.82"
$ws.Range("D22").Style = $normalStyle
$ws.Range("E22").Value = "  -6.17%  "

$ws.Range("E23").Value = "  -5.88%  "

$ws.Range("D24").Value = "'77.15"
$ws.Range("D24").Style = $normalStyle
$ws.Range("E24").Value = "  -2.84%  "

$ws.Range("E25").Value = "  -0.01%  "

$ws.Range("D26").Value = "'3.525.44"
$ws.Range("D26").Style = $normalStyle
$ws.Range("E26").Value = "  -2.40%  "

$ws.Range("E27").Value = "  -9.85%  "

$ws.Range("D28").Value = "'9.18"
$ws.Range("D28").Style = $normalStyle
$ws.Range("E28").Value = "  -6.43%  "

$ws.Range("D29").Value = "'7.74"
$ws.Range("D29").Style = $normalStyle
$ws.Range("E29").Value = "  -7.91%  "

$ws.Range("E30").Value = "  -3.25%  "

$ws.Range("E31").Value = "  +0.18%  "

$ws.Range("D32").Value = "'0.160"
$ws.Range("D32").Style = $normalStyle
$ws.Range("E32").Value = "  -4.85%  "

$ws.Range("D33").Value = "'1.45"
$ws.Range("D33").Style = $normalStyle
$ws.Range("E33").Value = "  -8.78%  "

$ws.Range("D34").Value = "'24.26"
$ws.Range("D34").Style = $normalStyle
$ws.Range("E34").Value = "  -4.74%  "

$ws.Range("D35").Value = "'3.383.88"
$ws.Range("D35").Style = $normalStyle
$ws.Range("E35").Value = "  -2.36%  "

$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "'1.67"
$ws.Range("D37").Style = $normalStyle
$ws.Range("E37").Value = "  -7.39%  "

$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D38").Value = "'5.50"
$ws.Range("D38").Style = $normalStyle
$ws.Range("E38").Value = "  -9.46%  "

$ws.Range("B39").Value = "FirstDigitalUSD"
$ws.Range("C39").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D39").Value = "'0.999"
$ws.Range("D39").Style = $normalStyle
$ws.Range("E39").Value = "  -0.03%  "

$ws.Range("B40").Value = "Aptos"
$ws.Range("C40").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D40").Value = "'7.47"
$ws.Range("D40").Style = $normalStyle
$ws.Range("E40").Value = "  -6.04%  "

$ws.Range("D41").Value = "'167.93"
$ws.Range("D41").Style = $normalStyle
$ws.Range("E41").Value = "  -4.99%  "

$ws.Range("D42").Value = "'0.0850"
$ws.Range("D42").Style = $normalStyle
$ws.Range("E42").Value = "  -4.61%  "

$ws.Range("D43").Value = "'0.867"
$ws.Range("D43").Style = $normalStyle
$ws.Range("E43").Value = "  -2.32%  "

$ws.Range("E44").Value = "  -8.00%  "

$ws.Range("D45").Value = "'1.91"
$ws.Range("D45").Style = $normalStyle
$ws.Range("E45").Value = "  -10.76%  "

$ws.Range("D46").Value = "'45.33"
$ws.Range("D46").Style = $normalStyle
$ws.Range("E46").Value = "  -2.13%  "

$ws.Range("D47").Value = "'26.38"
$ws.Range("D47").Style = $normalStyle
$ws.Range("E47").Value = "  -10.19%  "

$ws.Range("E48").Value = "  -5.44%  "

$ws.Range("D49").Value = "'7.01"
$ws.Range("D49").Style = $normalStyle
$ws.Range("E49").Value = "  -6.38%  "

$ws.Range("D50").Value = "'2.25"
$ws.Range("D50").Style = $normalStyle
$ws.Range("E50").Value = "  -8.81%  "

$ws.Range("D51").Value = "'0.912"
$ws.Range("D51").Style = $normalStyle
$ws.Range("E51").Value = "  -7.67%  "
